# Generate Report for Handback
# Adds a new "handback" row (file 9580c269-5f72-4390-af45-3e62e0c8fa04) to the
# Overview / zh-cn / de-de sheets, and refreshes the existing row's
# (1f6f0656-7907-4df0-973b-310923028b4d, formerly b964ff0c-27b0-4326-8a2b-cb625594757d)
# generated-file names / timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "b964ff0c-27b0-4326-8a2b-cb625594757d"
$renamedGuid = "1f6f0656-7907-4df0-973b-310923028b4d"
$newGuid = "9580c269-5f72-4390-af45-3e62e0c8fa04"

$oldXlfHash = "2252a6a2323aa890a8ceef73037cfa85d9e21fa2"
$renamedXlfHash = "ec751c8ba71c347c52a0e12f808c2e11ead3b8b9"
$newXlfHash = "9dfcdf2c01952b0a1640ca0aa23b8a8de3e8189d"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2: rename the existing handback file + bump the "Latest HO Xliff
# Generate Date" timestamp.
$wsOverview.Range("A2").Value = "$renamedGuid.md"
$wsOverview.Range("B2").Value = "e2e\$renamedGuid.md"
$wsOverview.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/$renamedGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$renamedGuid.md"
$wsOverview.Range("G2").Value = "2016-08-13 21:19:36"
$wsOverview.Range("G2").NumberFormat = $dateFmt

# Row 3: new handback file.
$wsOverview.Range("A3").Value = "$newGuid.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-13 21:19:36"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/$newGuid.md", "", "", "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 updates (rename + refreshed xliff/timestamps).
$wsZh.Range("A2").Value = "$renamedGuid.md"
$wsZh.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/$renamedGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$renamedGuid.md"
$wsZh.Range("G2").Value = "$renamedGuid.$renamedXlfHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-13 21:19:27"
$wsZh.Range("H2").NumberFormat = $dateFmt
$wsZh.Range("I2").Value = "$renamedGuid.md"
$wsZh.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3c0863da44151c9802fc748488aebd921b1a380d/e2e/$renamedGuid.md"
$wsZh.Hyperlinks.Item(2).TextToDisplay = "$renamedGuid.md"
$wsZh.Range("J2").Value = "$renamedGuid.$renamedXlfHash.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-13 21:19:55"
$wsZh.Range("K2").NumberFormat = $dateFmt

# Row 3: new duplicate-content handback file.
$wsZh.Range("A3").Value = "$newGuid.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newGuid.$newXlfHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-13 21:19:27"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("I3").Value = "$newGuid.md"
$wsZh.Range("J3").Value = "$newGuid.$newXlfHash.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-13 21:19:55"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/$newGuid.md", "", "", "$newGuid.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3c0863da44151c9802fc748488aebd921b1a380d/e2e/$newGuid.md", "", "", "$newGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 updates (rename + refreshed xliff/timestamps).
$wsDe.Range("A2").Value = "$renamedGuid.md"
$wsDe.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/$renamedGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$renamedGuid.md"
$wsDe.Range("G2").Value = "$renamedGuid.$renamedXlfHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-13 21:19:36"
$wsDe.Range("H2").NumberFormat = $dateFmt
$wsDe.Range("I2").Value = "$renamedGuid.md"
$wsDe.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e8ede0ee116909653671d070ec3eb53ef9cc0062/e2e/$renamedGuid.md"
$wsDe.Hyperlinks.Item(2).TextToDisplay = "$renamedGuid.md"
$wsDe.Range("J2").Value = "$renamedGuid.$renamedXlfHash.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-13 21:20:11"
$wsDe.Range("K2").NumberFormat = $dateFmt

# Row 3: new duplicate-content handback file.
$wsDe.Range("A3").Value = "$newGuid.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newGuid.$newXlfHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-13 21:19:36"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("I3").Value = "$newGuid.md"
$wsDe.Range("J3").Value = "$newGuid.$newXlfHash.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-13 21:20:11"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/$newGuid.md", "", "", "$newGuid.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e8ede0ee116909653671d070ec3eb53ef9cc0062/e2e/$newGuid.md", "", "", "$newGuid.md") | Out-Null
